# "Remove temporal dimension from sensitivity array"
#
# The workbook originally had two near-identical sheets holding the
# sensitivity-coefficient array, one per angular bin ("sensitivity coef,
# ang1" and "sensitivity coef, ang2"), each still carrying rows for a
# temporal axis (Static/LF1/LF2/LF3/MF/HF). The edit collapses this to a
# single "sensitivity coef" sheet whose rows are now indexed by angle
# (angle 1..angle 6) instead of by time bin, and updates the Glossary to
# describe the angular bins as a user-supplied array instead of a
# min/max-angle lookup table.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# 1. Drop the second angular-bin sheet entirely.
$wb.Worksheets.Item("sensitivity coef, ang2").Delete()

# 2. The remaining sensitivity sheet loses its "ang1" qualifier since
#    there is now only one.
$wsSens = $wb.Worksheets.Item("sensitivity coef, ang1")
$wsSens.Name = "sensitivity coef"

# 3. Re-label the rows: what used to be temporal bins (Static/LF1/LF2/
#    LF3/MF/HF) are now simply angle 1 .. angle 6. The data values
#    themselves (columns B:N) are untouched.
$wsSens.Range("A2").Value = "angle 1"
$wsSens.Range("A3").Value = "angle 2"
$wsSens.Range("A4").Value = "angle 3"
$wsSens.Range("A5").Value = "angle 4"
$wsSens.Range("A6").Value = "angle 5"
$wsSens.Range("A7").Value = "angle 6"

$wsSens.Activate()
$wsSens.Range("F16").Select()

# 4. Update the Glossary's "ANGULAR BINS" section: remove the old
#    Label/Min/Max lookup table (rows 18-20) and instead note that the
#    angular bins are a user-defined array, in mas units.
$wsGloss = $wb.Worksheets.Item("Glossary")
$wsGloss.Range("D17").Value = "User defined array in mas units"
$wsGloss.Range("D17").Font.Bold = $false
$wsGloss.Range("A18:C18").ClearContents()
$wsGloss.Range("A19:C20").ClearContents()

$wsGloss.Activate()
$wsGloss.Range("C21").Select()
